$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '28.048.75'
$c.Style = "Normal"
$ws.Range('E2').Value = '  -0.44%  '

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.900.74'
$c.Style = "Normal"
$ws.Range('E3').Value = '  +1.69%  '

$ws.Range('E4').Value = '  -0.02%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '312.40'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +0.24%  '

$ws.Range('E6').Value = '  -0.08%  '

$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.5063'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +0.05%  '

$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.3917'
$c.Style = "Normal"
$ws.Range('E8').Value = '  +0.05%  '

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.09250'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -4.15%  '

$ws.Range('E10').Value = '  -0.31%  '

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '41.85'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +2.47%  '

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '6.356'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -2.09%  '

$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '20.78'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -0.62%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '1.899.94'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +1.33%  '

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +0.00%  '

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '7.289'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -1.80%  '

$ws.Range('E17').Value = '  -0.84%  '

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '92.27'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -0.70%  '

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.06569'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -0.76%  '

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '17.74'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +1.16%  '

$ws.Range('E21').Value = '  +0.00%  '

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '6.208'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +0.84%  '

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '28.112.42'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -0.43%  '

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '11.35'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +0.00%  '

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.321'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +1.62%  '

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.591'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +2.27%  '

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '2.120.91'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +1.48%  '

$ws.Range('E28').Value = '  -1.20%  '

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '157.51'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -0.18%  '

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '127.00'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -0.14%  '

$ws.Range('E31').Value = '  +1.46%  '

$ws.Range('E32').Value = '  +0.97%  '

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '5.589'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -0.60%  '

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '3.608'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -0.46%  '

$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '9.578'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -0.19%  '

$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.06663'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -1.04%  '

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.02403'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +0.80%  '

$ws.Range('E38').Value = '  -1.03%  '

$ws.Range('E39').Value = '  -0.29%  '

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '1.250'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +6.36%  '

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.6354'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +0.15%  '

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '11.43'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -0.28%  '

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '4.968'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -0.09%  '

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -0.10%  '

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '13.22'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -2.48%  '

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.5966'
$c.Style = "Normal"

$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '3.703'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +1.16%  '

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.273'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +0.67%  '

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '2.006'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +0.71%  '

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '122.42'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -1.55%  '

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '1.176'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -1.55%  '
